{"js": "// Replace each old division-problem string with its corresponding new one.\n// Every old string is unique within the document body, so a direct\n// search+replace keyed on the exact text is safe and unambiguous.\nconst replacements = [\n  [\"33\u00f73=11, 0\", \"70\u00f72=35, 0\"],\n  [\"25\u00f73=8, 1\", \"10\u00f74=2, 2\"],\n  [\"95\u00f78=11, 7\", \"96\u00f75=19, 1\"],\n  [\"20\u00f77=2, 6\", \"32\u00f77=4, 4\"],\n  [\"10\u00f78=1, 2\", \"32\u00f78=4, 0\"],\n  [\"13\u00f74=3, 1\", \"93\u00f72=46, 1\"],\n  [\"11\u00f75=2, 1\", \"73\u00f78=9, 1\"],\n  [\"51\u00f72=25, 1\", \"30\u00f73=10, 0\"],\n  [\"52\u00f73=17, 1\", \"47\u00f75=9, 2\"],\n  [\"78\u00f76=13, 0\", \"15\u00f78=1, 7\"],\n  [\"21\u00f78=2, 5\", \"39\u00f79=4, 3\"],\n  [\"66\u00f74=16, 2\", \"28\u00f79=3, 1\"],\n  [\"34\u00f76=5, 4\", \"88\u00f79=9, 7\"],\n  [\"81\u00f75=16, 1\", \"87\u00f73=29, 0\"],\n  [\"76\u00f79=8, 4\", \"15\u00f73=5, 0\"],\n  [\"65\u00f72=32, 1\", \"18\u00f78=2, 2\"],\n  [\"40\u00f77=5, 5\", \"16\u00f74=4, 0\"],\n  [\"35\u00f74=8, 3\", \"34\u00f77=4, 6\"],\n  [\"10\u00f72=5, 0\", \"41\u00f79=4, 5\"],\n  [\"76\u00f76=12, 4\", \"74\u00f77=10, 4\"],\n  [\"51\u00f73=17, 0\", \"48\u00f78=6, 0\"],\n  [\"23\u00f79=2, 5\", \"31\u00f74=7, 3\"],\n  [\"57\u00f76=9, 3\", \"20\u00f75=4, 0\"],\n  [\"65\u00f78=8, 1\", \"19\u00f78=2, 3\"],\n  [\"67\u00f79=7, 4\", \"91\u00f77=13, 0\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each old division-problem string with its corresponding new one.\n# Every old string occurs exactly once in the document body, so a plain\n# Find/Replace keyed on the literal text is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = '33\u00f73=11, 0'; New = '70\u00f72=35, 0' }\n    @{ Old = '25\u00f73=8, 1'; New = '10\u00f74=2, 2' }\n    @{ Old = '95\u00f78=11, 7'; New = '96\u00f75=19, 1' }\n    @{ Old = '20\u00f77=2, 6'; New = '32\u00f77=4, 4' }\n    @{ Old = '10\u00f78=1, 2'; New = '32\u00f78=4, 0' }\n    @{ Old = '13\u00f74=3, 1'; New = '93\u00f72=46, 1' }\n    @{ Old = '11\u00f75=2, 1'; New = '73\u00f78=9, 1' }\n    @{ Old = '51\u00f72=25, 1'; New = '30\u00f73=10, 0' }\n    @{ Old = '52\u00f73=17, 1'; New = '47\u00f75=9, 2' }\n    @{ Old = '78\u00f76=13, 0'; New = '15\u00f78=1, 7' }\n    @{ Old = '21\u00f78=2, 5'; New = '39\u00f79=4, 3' }\n    @{ Old = '66\u00f74=16, 2'; New = '28\u00f79=3, 1' }\n    @{ Old = '34\u00f76=5, 4'; New = '88\u00f79=9, 7' }\n    @{ Old = '81\u00f75=16, 1'; New = '87\u00f73=29, 0' }\n    @{ Old = '76\u00f79=8, 4'; New = '15\u00f73=5, 0' }\n    @{ Old = '65\u00f72=32, 1'; New = '18\u00f78=2, 2' }\n    @{ Old = '40\u00f77=5, 5'; New = '16\u00f74=4, 0' }\n    @{ Old = '35\u00f74=8, 3'; New = '34\u00f77=4, 6' }\n    @{ Old = '10\u00f72=5, 0'; New = '41\u00f79=4, 5' }\n    @{ Old = '76\u00f76=12, 4'; New = '74\u00f77=10, 4' }\n    @{ Old = '51\u00f73=17, 0'; New = '48\u00f78=6, 0' }\n    @{ Old = '23\u00f79=2, 5'; New = '31\u00f74=7, 3' }\n    @{ Old = '57\u00f76=9, 3'; New = '20\u00f75=4, 0' }\n    @{ Old = '65\u00f78=8, 1'; New = '19\u00f78=2, 3' }\n    @{ Old = '67\u00f79=7, 4'; New = '91\u00f77=13, 0' }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 0  # wdFindStop - don't wrap, don't touch text outside range\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n"}
